# Prepared last run of 600 Wind Samples for 4200 time limit
#
# Rework the small summary block on "Tabelle2": the Mean/Std-dev/CI helper
# table that used to live in columns I/J/L/M/N (rows 1-8) moves over to
# columns L/M/O/P/Q, and the fixed z=1.96 normal-approximation 95% CI is
# replaced by a proper 0.05 t-distribution based CI driven by T.INV.
# A couple of stray formatted-but-empty cells are removed too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# -----------------------------------------------------------------
# Wipe the old helper block and the stray empty formatted cells
# (leave I1 itself alone for a moment - its bold/left style gets
# reused by the new header label below).
# -----------------------------------------------------------------
$ws.Range("I2:N8").Clear()
$ws.Range("I9").Clear()
$ws.Range("J11").Clear()
$ws.Range("I17").Clear()
$ws.Range("J19").Clear()
$ws.Rows(27).Delete()
$ws.Rows(25).Delete()

# -----------------------------------------------------------------
# New header label describing the CI block (re-use the bold/left
# style that used to live on I1 by copying its formatting over),
# then clear out the now-redundant I1 cell.
# -----------------------------------------------------------------
$ws.Range("I1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "0.05 t-distribution 95% CI intervals"
$ws.Range("I1").Clear()

# -----------------------------------------------------------------
# Row 2: mean labels/formula + lower/upper bound headers.
# -----------------------------------------------------------------
$ws.Range("L2").Value = "Mean Obj:"
$ws.Range("M2").Formula = "=SUM(E2:E21)/20"
$ws.Range("O2").Value = "untere Grenze:"
$ws.Range("Q2").Value = "obere:"

# -----------------------------------------------------------------
# Row 3: mean labels/formula + Obj CI (t-distribution based).
# -----------------------------------------------------------------
$ws.Range("L3").Value = "Mean Gap:"
$ws.Range("M3").Formula = "=SUM(D2:D21)/20"
$ws.Range("O3").Value = "Obj:"
$ws.Range("P3").Formula = "=M2+(-H6)*M6/(SQRT(20))"
$ws.Range("Q3").Formula = "=M2+H6*M6/(SQRT(20))"

# -----------------------------------------------------------------
# Row 4: mean labels/formula + Gap CI.
# -----------------------------------------------------------------
$ws.Range("L4").Value = "Mean Time:"
$ws.Range("M4").Formula = "=SUM(F2:F21)/20"
$ws.Range("O4").Value = "Gap:"
$ws.Range("P4").Formula = "=M3-H6*M7/(SQRT(20))"
$ws.Range("Q4").Formula = "=M3+H6*M7/(SQRT(20))"

# -----------------------------------------------------------------
# Row 5: std-dev label + Time CI.
# -----------------------------------------------------------------
$ws.Range("L5").Value = "Standardabweichung:"
$ws.Range("O5").Value = "Time:"
$ws.Range("P5").Formula = "=M4-H7*M8/(SQRT(20))"
$ws.Range("Q5").Formula = "=M4+H7*M8/(SQRT(20))"

# -----------------------------------------------------------------
# Row 6: t critical value for the 95% two-sided CI + its confidence
# level, plus Obj std-dev label/formula.
# -----------------------------------------------------------------
$ws.Range("H6").Formula = "=T.INV(0.975,19)"
$ws.Range("I6").Value = 0.95
$ws.Range("I6").NumberFormat = "0%"
$ws.Range("L6").Value = "Obj:"
$ws.Range("M6").Formula = "=STDEV.P(E2:E21)"

# -----------------------------------------------------------------
# Row 7: t critical value for the 97.5% one-sided cut (Gap/Time CIs)
# plus Gap std-dev label/formula.
# -----------------------------------------------------------------
$ws.Range("H7").Formula = "=T.INV(0.9875,19)"
$ws.Range("I7").Value = 0.975
$ws.Range("I7").NumberFormat = "0%"
$ws.Range("L7").Value = "Gap:"
$ws.Range("M7").Formula = "=STDEV.P(D2:D21)"

# -----------------------------------------------------------------
# Row 8: Time std-dev label/formula.
# -----------------------------------------------------------------
$ws.Range("L8").Value = "Time:"
$ws.Range("M8").Formula = "=STDEV.P(F2:F21)"

# -----------------------------------------------------------------
# Selection / view bookkeeping to match the saved workbook state.
# -----------------------------------------------------------------
$ws.Range("H1:Q1048576").Select()
